$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (value does not look like a pure number, safe to assign directly)
$textUpdates = @{
    'D2' = '26.705.92'
    'E2' = '  -8.13%  '
    'D3' = '1.645.26'
    'E3' = '  -9.59%  '
    'E4' = '  +1.53%  '
    'E5' = '  -5.62%  '
    'E6' = '  +1.60%  '
    'E7' = '  -15.87%  '
    'E8' = '  -6.28%  '
    'E9' = '  -5.43%  '
    'E10' = '  -9.82%  '
    'E11' = '  -1.14%  '
    'D12' = '1.651.90'
    'E12' = '  -9.00%  '
    'E13' = '  -4.58%  '
    'E14' = '  -7.90%  '
    'D15' = '1.867.73'
    'E15' = '  -9.54%  '
    'E16' = '  -12.46%  '
    'E17' = '  -13.07%  '
    'D18' = '26.747.57'
    'E18' = '  -7.10%  '
    'E19' = '  -7.27%  '
    'E20' = '  +1.65%  '
    'E21' = '  -5.00%  '
    'E22' = '  -10.78%  '
    'E23' = '  +1.63%  '
    'E24' = '  -8.45%  '
    'E25' = '  -7.11%  '
    'E26' = '  -2.58%  '
    'E27' = '  -9.02%  '
    'E28' = '  -6.42%  '
    'E29' = '  -4.52%  '
    'E30' = '  -8.82%  '
    'E31' = '  -6.49%  '
    'E32' = '  -6.37%  '
    'E33' = '  -8.06%  '
    'E34' = '  -7.22%  '
    'E35' = '  -6.78%  '
    'E36' = '  -3.05%  '
    'E37' = '  -5.81%  '
    'E38' = '  -5.11%  '
    'E39' = '  -7.37%  '
    'B40' = 'PaxDollar'
    'C40' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'E40' = '  +1.73%  '
    'B41' = 'TrustWalletToken'
    'C41' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E41' = '  -0.79%  '
    'D42' = '1.059.48'
    'E42' = '  -6.20%  '
    'E43' = '  -10.03%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E44' = '  -4.05%  '
    'B45' = 'RocketPoolETH'
    'C45' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D45' = '1.785.80'
    'E45' = '  -9.46%  '
    'E46' = '  -0.87%  '
    'E47' = '  -2.35%  '
    'B48' = 'Frax'
    'C48' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'E48' = '  -1.35%  '
    'B49' = 'Aave'
    'C49' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'E49' = '  -8.29%  '
    'B50' = 'Cronos'
    'C50' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'E50' = '  -4.36%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E51' = '  -5.32%  '
}

# Numeric-looking text updates (must be forced to stay text, matching original inlineStr type)
$numericTextUpdates = @{
    'D5' = '219.68'
    'D6' = '1.022'
    'D7' = '0.4954'
    'D8' = '0.2557'
    'D9' = '21.71'
    'D10' = '0.06113'
    'D11' = '0.07411'
    'D13' = '4.426'
    'D14' = '0.5714'
    'D16' = '0.000008031'
    'D17' = '64.34'
    'D19' = '4.997'
    'D20' = '1.023'
    'D21' = '10.76'
    'D22' = '184.13'
    'D23' = '1.023'
    'D24' = '6.177'
    'D25' = '142.98'
    'D26' = '7.530'
    'D27' = '0.1145'
    'D28' = '15.10'
    'D29' = '1.345'
    'D30' = '0.05757'
    'D31' = '1.334'
    'D32' = '3.416'
    'D33' = '3.392'
    'D34' = '1.551'
    'D35' = '0.9689'
    'D36' = '2.455'
    'D37' = '0.5913'
    'D38' = '2.608'
    'D39' = '0.01567'
    'D40' = '1.024'
    'D41' = '0.8586'
    'D43' = '5.751'
    'D44' = '95.40'
    'D46' = '0.00000000111'
    'D47' = '0.4421'
    'D48' = '0.9993'
    'D49' = '54.74'
    'D50' = '0.05244'
    'D51' = '7.709'
}

foreach ($key in $textUpdates.Keys) {
    $ws.Range($key).Value2 = $textUpdates[$key]
}

foreach ($key in $numericTextUpdates.Keys) {
    $ws.Range($key).Value2 = "'" + $numericTextUpdates[$key]
    $ws.Range($key).Style = "Normal"
}
